# Apply updated cryptocurrency price/volume figures (and the
# MultiversX/Celestia rank swap) to the worksheet, cell by cell,
# mirroring the upstream data-refresh commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.215.27"
$ws.Range("E2").Value = "  -1.59%  "
$ws.Range("D3").Value = "2.177.18"
$ws.Range("E3").Value = "  -2.33%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.26%  "
$ws.Range("D5").Value = "'250.82"
$ws.Range("E5").Value = "  +0.35%  "
$ws.Range("D6").Value = "'0.612"
$ws.Range("D7").Value = "'66.45"
$ws.Range("E7").Value = "  -7.80%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("D9").Value = "'0.576"
$ws.Range("E9").Value = "  -2.80%  "
$ws.Range("D10").Value = "'59.13"
$ws.Range("E10").Value = "  +1.77%  "
$ws.Range("D11").Value = "'36.39"
$ws.Range("E11").Value = "  -11.70%  "
$ws.Range("D12").Value = "'0.0936"
$ws.Range("E12").Value = "  -3.48%  "
$ws.Range("D13").Value = "'0.104"
$ws.Range("E13").Value = "  -1.26%  "
$ws.Range("D14").Value = "'6.87"
$ws.Range("E14").Value = "  -3.87%  "
$ws.Range("D15").Value = "2.501.29"
$ws.Range("E15").Value = "  -2.33%  "
$ws.Range("D16").Value = "'14.33"
$ws.Range("E16").Value = "  -4.41%  "
$ws.Range("D17").Value = "'0.848"
$ws.Range("E17").Value = "  -2.04%  "
$ws.Range("D18").Value = "2.177.02"
$ws.Range("E18").Value = "  -2.40%  "
$ws.Range("D19").Value = "41.070.28"
$ws.Range("E19").Value = "  -1.93%  "
$ws.Range("D20").Value = "0.0₃0951"
$ws.Range("E20").Value = "  -1.65%  "
$ws.Range("E21").Value = "  -1.66%  "
$ws.Range("D22").Value = "'6.08"
$ws.Range("E22").Value = "  -2.09%  "
$ws.Range("D23").Value = "'230.28"
$ws.Range("E23").Value = "  -2.20%  "
$ws.Range("E25").Value = "  -7.45%  "
$ws.Range("D26").Value = "'11.53"
$ws.Range("E26").Value = "  +7.92%  "
$ws.Range("D27").Value = "'1.00"
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("D28").Value = "'2.41"
$ws.Range("E28").Value = "  -4.62%  "
$ws.Range("D29").Value = "'168.12"
$ws.Range("E29").Value = "  -1.89%  "
$ws.Range("E30").Value = "  -3.49%  "
$ws.Range("D31").Value = "'20.25"
$ws.Range("E31").Value = "  -2.43%  "
$ws.Range("D32").Value = "'0.121"
$ws.Range("E32").Value = "  -3.03%  "
$ws.Range("D33").Value = "'5.85"
$ws.Range("E33").Value = "  +5.23%  "
$ws.Range("E34").Value = "  +3.79%  "
$ws.Range("E35").Value = "  -3.23%  "
$ws.Range("D36").Value = "'4.53"
$ws.Range("E36").Value = "  -3.80%  "
$ws.Range("E37").Value = "  +0.25%  "
$ws.Range("D38").Value = "'24.68"
$ws.Range("E38").Value = "  -5.86%  "
$ws.Range("E39").Value = "  +0.89%  "
$ws.Range("D40").Value = "'2.22"
$ws.Range("E40").Value = "  -3.00%  "
$ws.Range("D41").Value = "'5.29"
$ws.Range("E41").Value = "  +6.87%  "
$ws.Range("D42").Value = "'5.49"
$ws.Range("E42").Value = "  -8.56%  "
$ws.Range("B43").Value = "Celestia"
$ws.Range("C43").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D43").Value = "'11.44"
$ws.Range("E43").Value = "  -4.51%  "
$ws.Range("B44").Value = "MultiversX"
$ws.Range("C44").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D44").Value = "'60.87"
$ws.Range("E44").Value = "  -9.51%  "
$ws.Range("D45").Value = "'8.51"
$ws.Range("E45").Value = "  -3.05%  "
$ws.Range("D46").Value = "'0.1000"
$ws.Range("E46").Value = "  -1.56%  "
$ws.Range("E47").Value = "  -6.45%  "
$ws.Range("E48").Value = "  -0.29%  "
$ws.Range("E49").Value = "  -1.76%  "
$ws.Range("E50").Value = "  -4.74%  "
$ws.Range("D51").Value = "'2.74"
$ws.Range("E51").Value = "  +1.19%  "

# Strip the "quote prefix" style Excel applies when a text value
# looks numeric, so the cells' formatting stays untouched (General).
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D51").Style = "Normal"
